# Fix for crash in the Kenya database: cast referral fid as a character.
# Adds two new rows to the day7 dictionary sheet, right after the
# "a1-contact-success" entry, describing two new fields:
#   a1-contact-exist    -> valid_phone -> exist
#   a1-contact-phoneoff -> phone_off   -> phoneoff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 17 (pushing everything from the old
# row 17 down to row 19). Inserting picks up formatting from the row
# directly above, matching the existing "4 / 1 / 4" style pattern used
# by the surrounding dictionary rows.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# New row 17: a1-contact-exist / valid_phone / exist
$ws.Range("A17").Value = "a1-contact-exist"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "valid_phone"
$ws.Range("D17").Value = "exist"

# New row 18: a1-contact-phoneoff / phone_off / phoneoff
$ws.Range("A18").Value = "a1-contact-phoneoff"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "phone_off"
$ws.Range("D18").Value = "phoneoff"

# Match the author's last-saved cursor position.
$ws.Range("E18").Select()
